# merge the errortypes as a alias of Error.xlsx
# Insert a new "别名" (Alias) column between the existing "错误id" and
# "描述" columns of table "表1" on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- 1. Grow the table by one column -------------------------------------
# ListColumns.Add() appends a new (blank) column at the right-hand edge of
# the table, which becomes worksheet column C here (table was A:B, 26 rows).
$newCol = $lo.ListColumns.Add()

# --- 2. Shift the existing "描述" column (B) into the new column (C) -----
# Range.Copy carries both the value and the cell style/format along, so the
# header-row formatting (rows 2 & 3) follows the data into column C.
$ws.Range("B1:B26").Copy($ws.Range("C1:C26"))
# Touch C1 again directly so the table's column-3 name re-syncs from the
# header cell text ("描述") instead of staying "Column3".
$ws.Cells.Item(1, 3).Value2 = $ws.Cells.Item(1, 3).Value2

# --- 3. Fill column B with the new "别名" (alias) values ------------------
$aliasValues = @(
    "别名",
    "string",
    "Alias",
    "OK",
    "CommonError",
    "BattleNoUseCard",
    "BattleNoUseSpellCard",
    "BattleLackMp",
    "BattleLackLp",
    "BattleLackPp",
    "BattleHeroSkillInCd",
    "CardOutPunish",
    "CardFullPunish",
    "DeckCardTypeLimitLegend",
    "DeckCardTypeLimit",
    "DeckIsFull",
    "CardExpNotEnough",
    "CardExpNotEnough2",
    "CardJobTwice",
    "SceneLevelNeed",
    "SceneAPNotEnough",
    "SceneWarpNeedActive",
    "BagNotEnoughDimond",
    "BagNotEnoughResource",
    "BagIsFull",
    "BagNotEnoughItems"
)
for ($r = 1; $r -le 26; $r++) {
    $ws.Cells.Item($r, 2).Value2 = $aliasValues[$r - 1]
}

# --- 4. Column widths (B narrower, C widened to fit the old "描述" text) --
$ws.Columns.Item(2).ColumnWidth = 14
$ws.Columns.Item(3).ColumnWidth = 27

# --- 5. Restore the workbook's remembered selection -----------------------
$ws.Range("C10").Select()
